$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns stay formatted as text, matching the
# original inline-string cell type, before assigning new values so that
# Excel does not reinterpret numeric-looking strings (e.g. "312.39") as
# numbers or dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.484.81"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "2.471.70"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").Value = "312.39"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "91.49"
$ws.Range("E6").Value = "  -7.09%  "
$ws.Range("D7").Value = "0.541"
$ws.Range("E7").Value = "  -3.69%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -4.90%  "
$ws.Range("D10").Value = "32.76"
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "2.846.93"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "6.82"
$ws.Range("E14").Value = "  -5.38%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.514.08"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.19"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").Value = "41.333.97"
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  -4.85%  "
$ws.Range("D20").Value = "0.0₃0916"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").Value = "70.43"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("E22").Value = "  -9.88%  "
$ws.Range("D23").Value = "234.41"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  -6.22%  "
$ws.Range("D27").Value = "23.88"
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "9.66"
$ws.Range("D30").Value = "35.95"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").Value = "152.48"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "5.40"
$ws.Range("E32").Value = "  -8.49%  "
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "2.54"
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("D35").Value = "0.0751"
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "2.98"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("E38").Value = "  -6.96%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("D40").Value = "0.0990"
$ws.Range("E40").Value = "  -8.22%  "
$ws.Range("D41").Value = "4.03"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "19.24"
$ws.Range("E43").Value = "  -11.29%  "
$ws.Range("D44").Value = "1.957.65"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "0.0281"
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  -8.77%  "
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "2.707.42"
$ws.Range("D49").Value = "95.31"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("D50").Value = "67.59"
$ws.Range("E50").Value = "  -6.02%  "
$ws.Range("E51").Value = "  -6.98%  "
